# "Upload Baitap1 lan 2" - re-upload of the grading sheet with cleaned-up
# answer text (the previous upload stored Python byte-string reprs like
# b'explore' instead of plain text), one grade correction (L5, "monkey"
# question) flipped from correct to wrong, and the resulting score in L12
# updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ----- Row 2 : "explore" -----
$ws.Range("B2:K2").Value = "explore"
$ws.Range("L2").Value = '\xe0\xb4\x86\xe0\xb4\xb0'

# ----- Row 3 : "hello" -----
$ws.Range("B3:L3").Value = "hello"

# ----- Row 4 : "joyful" (J4 / L4 are distinct answers) -----
$ws.Range("B4:I4").Value = "joyful"
$ws.Range("J4").Value = "mt hol"
$ws.Range("K4").Value = "joyful"
$ws.Range("L4").Value = "my holi"

# ----- Row 5 : "monkey" (L5 is the mis-typed "Monkey " answer) -----
$ws.Range("B5:K5").Value = "monkey"
$ws.Range("L5").Value = "Monkey "

# ----- Row 6 : "pig" (K6 / L6 are the "bid" answer) -----
$ws.Range("B6:J6").Value = "pig"
$ws.Range("K6:L6").Value = "bid"

# ----- Row 7 : "question" (L7 is "Question ") -----
$ws.Range("B7:K7").Value = "question"
$ws.Range("L7").Value = "Question "

# ----- Row 8 : "vehicle" -----
$ws.Range("B8:L8").Value = "vehicle"

# ----- Row 9 : "word" (L9 is "Word ") -----
$ws.Range("B9:K9").Value = "word"
$ws.Range("L9").Value = "Word "

# ----- Row 10 : "yatch" -----
$ws.Range("B10:L10").Value = "yatch"

# ----- Row 11 : "zebra" -----
$ws.Range("B11:I11").Value = "zebra"

# ----- Row 12 label -----
$ws.Range("B12").Value = "% Correct"

# L5 ("Monkey ") is now graded as wrong -> highlight with a solid red fill,
# matching the other grading colors already used by the conditional
# formatting (dxf) rules on this sheet.
$ws.Range("L5").Interior.Color = 255

# The % Correct score for that column drops accordingly.
$ws.Range("L12").Value = 0.3

# Selection left on the grid like the re-uploaded workbook.
$ws.Range("B1:L12").Select()
